$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "www.test"
$ws.Range("B7").Value = "lena_dima"
$ws.Range("C7").Value = "&).A-O}\"
